$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the worker detail table (rows 16-21) so that the two records
# belonging to period 1710 (FERNANDO FRANCISCO FELFLE FUENTES and
# CARLOS ALBERTO HOYOS RIOS) are listed first, followed by the period
# 1711 records (GUSTAVO ADOLFO FELFLE FUENTES, FERNANDO FRANCISCO
# FELFLE FUENTES, CARLOS ALBERTO HOYOS RIOS, OSCAR ENRIQUE ANAYA MEJIA).

$ws.Range("C16").Value2 = "9096389"
$ws.Range("D16").Value2 = "FERNANDO FRANCISCO FELFLE FUENTES"
$ws.Range("E16").Value2 = "1710"

$ws.Range("C17").Value2 = "9100677"
$ws.Range("D17").Value2 = "CARLOS ALBERTO HOYOS RIOS"
$ws.Range("E17").Value2 = "1710"

$ws.Range("C18").Value2 = "73199947"
$ws.Range("D18").Value2 = "GUSTAVO ADOLFO FELFLE FUENTES"
$ws.Range("E18").Value2 = "1711"

$ws.Range("C19").Value2 = "9096389"
$ws.Range("D19").Value2 = "FERNANDO FRANCISCO FELFLE FUENTES"
$ws.Range("E19").Value2 = "1711"

$ws.Range("C20").Value2 = "9100677"
$ws.Range("D20").Value2 = "CARLOS ALBERTO HOYOS RIOS"
$ws.Range("E20").Value2 = "1711"

$ws.Range("C21").Value2 = "1143393504"
$ws.Range("D21").Value2 = "OSCAR ENRIQUE ANAYA MEJIA"
$ws.Range("E21").Value2 = "1711"

$wb.Save()
